$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("D12").Value = 45049
$ws.Range("H12").Value = 'Madrigal'
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 17000
$ws.Range("L12").Value = 18000
$ws.Range("M12").Value = 17500
$ws.Range("N12").Value = '$/caja 40 unidades'
$ws.Range("O12").Value = 'Provincia del Elquí'
$ws.Range("P12").Value = 438
$ws.Range("Q12").Value = 40

# Row 13
$ws.Range("D13").Value = 44701
$ws.Range("H13").Value = 'Española'
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 19000
$ws.Range("L13").Value = 20000
$ws.Range("M13").Value = 19500
$ws.Range("N13").Value = '$/caja 30 unidades'
$ws.Range("O13").Value = 'Provincia del Elquí'
$ws.Range("P13").Value = 650
$ws.Range("Q13").Value = 30

# Row 14
$ws.Range("D14").Value = 44784
$ws.Range("H14").Value = 'Madrigal'
$ws.Range("J14").Value = 520
$ws.Range("K14").Value = 11500
$ws.Range("L14").Value = 12000
$ws.Range("M14").Value = 11750
$ws.Range("N14").Value = '$/caja 40 unidades'
$ws.Range("O14").Value = 'Provincia del Elquí'
$ws.Range("P14").Value = 294
$ws.Range("Q14").Value = 40

# Row 15
$ws.Range("D15").Value = 44839
$ws.Range("H15").Value = 'Española'
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 12000
$ws.Range("L15").Value = 13000
$ws.Range("M15").Value = 12500
$ws.Range("N15").Value = '$/caja 30 unidades'
$ws.Range("O15").Value = 'Provincia del Elquí'
$ws.Range("P15").Value = 417
$ws.Range("Q15").Value = 30

# Row 16
$ws.Range("D16").Value = 44420
$ws.Range("H16").Value = 'Madrigal'
$ws.Range("J16").Value = 800
$ws.Range("K16").Value = 14000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 14500
$ws.Range("N16").Value = '$/caja 40 unidades'
$ws.Range("O16").Value = 'Provincia de Limarí'
$ws.Range("P16").Value = 362
$ws.Range("Q16").Value = 40

# Row 17
$ws.Range("D17").Value = 44420
$ws.Range("H17").Value = 'Madrigal'
$ws.Range("J17").Value = 700
$ws.Range("K17").Value = 13000
$ws.Range("L17").Value = 14000
$ws.Range("M17").Value = 13500
$ws.Range("N17").Value = '$/caja 40 unidades'
$ws.Range("O17").Value = 'Provincia del Elquí'
$ws.Range("P17").Value = 338
$ws.Range("Q17").Value = 40

# Row 18
$ws.Range("D18").Value = 44687
$ws.Range("H18").Value = 'Española'
$ws.Range("J18").Value = 400
$ws.Range("K18").Value = 18000
$ws.Range("L18").Value = 19000
$ws.Range("M18").Value = 18500
$ws.Range("N18").Value = '$/caja 30 unidades'
$ws.Range("O18").Value = 'Provincia de Limarí'
$ws.Range("P18").Value = 617
$ws.Range("Q18").Value = 30

# Row 19
$ws.Range("D19").Value = 44858
$ws.Range("H19").Value = 'Española'
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 9500
$ws.Range("L19").Value = 10000
$ws.Range("M19").Value = 9750
$ws.Range("N19").Value = '$/caja 30 unidades'
$ws.Range("O19").Value = 'Provincia del Elquí'
$ws.Range("P19").Value = 325
$ws.Range("Q19").Value = 30

# Row 20
$ws.Range("D20").Value = 44498
$ws.Range("H20").Value = 'Española'
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 8500
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = 8750
$ws.Range("N20").Value = '$/caja 30 unidades'
$ws.Range("O20").Value = 'Provincia de Limarí'
$ws.Range("P20").Value = 292
$ws.Range("Q20").Value = 30

# Row 21
$ws.Range("A21").Value = 8
$ws.Range("B21").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C21").Value = 'Coquimbo'
$ws.Range("D21").Value = 44427
$ws.Range("D21").NumberFormat = $ws.Range("D20").NumberFormat
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 100112013
$ws.Range("G21").Value = 'Alcachofa'
$ws.Range("H21").Value = 'Madrigal'
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 400
$ws.Range("K21").Value = 12000
$ws.Range("L21").Value = 13000
$ws.Range("M21").Value = 12500
$ws.Range("N21").Value = '$/caja 40 unidades'
$ws.Range("O21").Value = 'Provincia de Limarí'
$ws.Range("P21").Value = 312
$ws.Range("Q21").Value = 40
$ws.Range("R21").Value = 'Hortaliza'
